$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 16 - this pushes the existing rows 16-23
# down to rows 17-24, matching the target dimension A1:R24.
$ws.Rows.Item(16).Insert()

# Populate the newly inserted row 16 with the new weekly price record.
$ws.Range("A16").Value = 3
$ws.Range("B16").Value = "Femacal de La Calera"
$ws.Range("C16").Value = "Coquimbo"
$ws.Range("D16").Value = 44468
$ws.Range("E16").Value = 5
$ws.Range("F16").Value = 100112022
$ws.Range("G16").Value = "Arveja Verde"
$ws.Range("H16").Value = "Perfection"
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 65
$ws.Range("K16").Value = 24000
$ws.Range("L16").Value = 25000
$ws.Range("M16").Value = 24538
$ws.Range("N16").Value = "$/malla 25 kilos"
$ws.Range("O16").Value = "Provincia de Limarí"
$ws.Range("P16").Value = 982
$ws.Range("Q16").Value = 25
$ws.Range("R16").Value = "Hortaliza"

# Match the original date cell number format (yyyy-mm-dd hh:mm:ss style)
# used by the other rows in column D.
$ws.Range("D16").NumberFormat = $ws.Range("D17").NumberFormat
